$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("B3").Value = 0.462043795620438
$ws.Range("C3").Value = 0.4854014598540146
$ws.Range("D3").Value = 0.2554744525547445
$ws.Range("E3").Value = 0.5182481751824818
$ws.Range("F3").Value = 0.1021897810218978
$ws.Range("G3").Value = 0.7467669631750652
$ws.Range("H3").Value = 0.7336931200496971
$ws.Range("I3").Value = 0.7325086438724548
$ws.Range("J3").Value = 0.4254385832780359
$ws.Range("K3").Value = 0.4333354317807579
$ws.Range("L3").Value = 0.6969385472475346
$ws.Range("M3").Value = 0.6752470491091697
$ws.Range("N3").Value = 0.2049129702414374
$ws.Range("O3").Value = 0.01267405675935353
$ws.Range("P3").Value = 0.1437286223068514
$ws.Range("Q3").Value = 0.1774106799846331
$ws.Range("R3").Value = 0.6724588735156335
$ws.Range("S3").Value = 0.6050479783482325
$ws.Range("T3").Value = 0.6103404989171414
$ws.Range("U3").Value = 0.5431357994144088
$ws.Range("V3").Value = 0.353697238971755
$ws.Range("W3").Value = 0.3949921463549847
$ws.Range("X3").Value = 0.4659854014598541
$ws.Range("Y3").Value = 0.6385401459854015
$ws.Range("Z3").Value = 0.6387222663872226
$ws.Range("AA3").Value = 0.2338071456012293
$ws.Range("AB3").Value = 0.2047217153284671
$ws.Range("AC3").Value = 0.6689925721084279
$ws.Range("AD3").Value = 0.6545534287772939
$ws.Range("AE3").Value = 0.7540233928414387
$ws.Range("AF3").Value = 0.7396007387213086
$ws.Range("AG3").Value = 0.745686794956868
$ws.Range("AH3").Value = 0.7205540809555409
$ws.Range("AI3").Value = 0.7423357664233576
$ws.Range("AJ3").Value = 0.7524259338772005
$ws.Range("AK3").Value = 0.4014598540145985
$ws.Range("AL3").Value = 0.458029197080292
$ws.Range("AM3").Value = 0.4233576642335766
$ws.Range("AN3").Value = 0.364963503649635
$ws.Range("AO3").Value = 0.3795620437956204
$ws.Range("AQ3").Value = 0.1313868613138686
$ws.Range("AR3").Value = 0.2043795620437956
$ws.Range("AS3").Value = 0.2043795620437956
$ws.Range("AT3").Value = 0.072992700729927
$ws.Range("AU3").Value = 0.072992700729927
$ws.Range("AV3").Value = 0.072992700729927
$ws.Range("AW3").Value = 0.08759124087591241
$ws.Range("AX3").Value = 0.05109489051094891
$ws.Range("AY3").Value = 0.06569343065693431
$ws.Range("AZ3").Value = 0.05109489051094891
$ws.Range("BA3").Value = 0.05109489051094891
$ws.Range("BB3").Value = 0.04379562043795621
$ws.Range("BC3").Value = 0.04379562043795621
$ws.Range("BD3").Value = 0.04379562043795621
$ws.Range("BE3").Value = 0.0583941605839416
$ws.Range("BF3").Value = 0.0583941605839416
$ws.Range("BG3").Value = 0.0583941605839416
$ws.Range("BH3").Value = 0.0583941605839416
$ws.Range("BI3").Value = 0.0583941605839416
$ws.Range("BJ3").Value = 0.06569343065693431
$ws.Range("BK3").Value = 0.06569343065693431
$ws.Range("BL3").Value = 0.06569343065693431
$ws.Range("BM3").Value = 0.06569343065693431
$ws.Range("BN3").Value = 0.06569343065693431
$ws.Range("BO3").Value = 0.06569343065693431
$ws.Range("BT3").Value = 0.05109489051094891

# Row 4
$ws.Range("B4").Value = 0.3461256811182618
$ws.Range("C4").Value = 0.4997852694768909
$ws.Range("D4").Value = 0.2689326935581176
$ws.Range("E4").Value = 0.5015005392792523
$ws.Range("F4").Value = 0.3040092764849821
$ws.Range("G4").Value = 0.2965072301990125
$ws.Range("H4").Value = 0.374112627417894
$ws.Range("I4").Value = 0.3737938413659845
$ws.Range("J4").Value = 0.1687738691158554
$ws.Range("K4").Value = 0.1810107455675833
$ws.Range("L4").Value = 0.2218866884584537
$ws.Range("M4").Value = 0.2386735989225567
$ws.Range("N4").Value = 0.1683819789862626
$ws.Range("O4").Value = 0.08503786694531396
$ws.Range("P4").Value = 0.1445811207769394
$ws.Range("Q4").Value = 0.2285815737140521
$ws.Range("R4").Value = 0.1968223587002977
$ws.Range("S4").Value = 0.2032568703546531
$ws.Range("T4").Value = 0.210620122106326
$ws.Range("U4").Value = 0.1894416995281113
$ws.Range("V4").Value = 0.2249185825300107
$ws.Range("W4").Value = 0.2076489890720712
$ws.Range("X4").Value = 0.245448149112929
$ws.Range("Y4").Value = 0.2386251679997726
$ws.Range("Z4").Value = 0.2422653732954648
$ws.Range("AA4").Value = 0.2130094680917154
$ws.Range("AB4").Value = 0.16711194210222
$ws.Range("AC4").Value = 0.2040565692171277
$ws.Range("AD4").Value = 0.1978498581907108
$ws.Range("AE4").Value = 0.2217953502494493
$ws.Range("AF4").Value = 0.2180417856311708
$ws.Range("AG4").Value = 0.1978062577872582
$ws.Range("AH4").Value = 0.2064010581405945
$ws.Range("AI4").Value = 0.201141321674305
$ws.Range("AJ4").Value = 0.2043834144943563
$ws.Range("AK4").Value = 0.4919925575506937
$ws.Range("AL4").Value = 0.2342493770603732
$ws.Range("AM4").Value = 0.4959042207278343
$ws.Range("AN4").Value = 0.4831866100471024
$ws.Range("AO4").Value = 0.4870588027869559
$ws.Range("AQ4").Value = 0.3390626839533289
$ws.Range("AR4").Value = 0.4047273252629636
$ws.Range("AS4").Value = 0.4047273252629636
$ws.Range("AT4").Value = 0.2610791094992419
$ws.Range("AU4").Value = 0.2610791094992419
$ws.Range("AV4").Value = 0.2610791094992419
$ws.Range("AW4").Value = 0.283736947218425
$ws.Range("AX4").Value = 0.2209993306887549
$ws.Range("AY4").Value = 0.2486546022226944
$ws.Range("AZ4").Value = 0.2209993306887549
$ws.Range("BA4").Value = 0.2209993306887549
$ws.Range("BB4").Value = 0.2053910595269901
$ws.Range("BC4").Value = 0.2053910595269901
$ws.Range("BD4").Value = 0.2053910595269901
$ws.Range("BE4").Value = 0.2353477826306845
$ws.Range("BF4").Value = 0.2353477826306845
$ws.Range("BG4").Value = 0.2353477826306845
$ws.Range("BH4").Value = 0.2353477826306845
$ws.Range("BI4").Value = 0.2353477826306845
$ws.Range("BJ4").Value = 0.2486546022226944
$ws.Range("BK4").Value = 0.2486546022226944
$ws.Range("BL4").Value = 0.2486546022226944
$ws.Range("BM4").Value = 0.2486546022226944
$ws.Range("BN4").Value = 0.2486546022226944
$ws.Range("BO4").Value = 0.2486546022226944
$ws.Range("BT4").Value = 0.2209993306887549

# Row 6
$ws.Range("B6").Value = 0.1545454545454545
$ws.Range("G6").Value = 0.7942122186495176
$ws.Range("H6").Value = 0.8829787234042553
$ws.Range("I6").Value = 0.8736842105263157
$ws.Range("J6").Value = 0.3513636363636364
$ws.Range("K6").Value = 0.3529019173996812
$ws.Range("L6").Value = 0.6570995701056114
$ws.Range("M6").Value = 0.628844839371155
$ws.Range("N6").Value = 0.1
$ws.Range("O6").Value = 0.003536345776031434
$ws.Range("P6").Value = 0.06763285024154589
$ws.Range("Q6").Value = 0.01052631578947368
$ws.Range("R6").Value = 0.6268656716417911
$ws.Range("S6").Value = 0.5561797752808989
$ws.Range("T6").Value = 0.5274725274725275
$ws.Range("U6").Value = 0.4632768361581921
$ws.Range("V6").Value = 0.2065217391304348
$ws.Range("W6").Value = 0.270042194092827
$ws.Range("X6").Value = 0.32
$ws.Range("Y6").Value = 0.5600000000000001
$ws.Range("Z6").Value = 0.5643564356435643
$ws.Range("AA6").Value = 0.0736842105263158
$ws.Range("AB6").Value = 0.078125
$ws.Range("AD6").Value = 0.6223175965665236
$ws.Range("AE6").Value = 0.7349397590361446
$ws.Range("AF6").Value = 0.7228915662650602
$ws.Range("AG6").Value = 0.7159090909090909
$ws.Range("AH6").Value = 0.6818181818181819
$ws.Range("AI6").Value = 0.7000000000000001
$ws.Range("AJ6").Value = 0.7176470588235294

# Row 7
$ws.Range("B7").Value = 0.3818181818181818
$ws.Range("C7").Value = 0
$ws.Range("G7").Value = 0.842443729903537
$ws.Range("H7").Value = 0.9095744680851063
$ws.Range("I7").Value = 0.9052631578947368
$ws.Range("J7").Value = 0.4357
$ws.Range("K7").Value = 0.4463838684108614
$ws.Range("L7").Value = 0.738456609705134
$ws.Range("M7").Value = 0.7328694463431304
$ws.Range("N7").Value = 0.1576923076923077
$ws.Range("O7").Value = 0.004911591355599214
$ws.Range("P7").Value = 0.1111111111111111
$ws.Range("Q7").Value = 0.08421052631578947
$ws.Range("R7").Value = 0.7164179104477612
$ws.Range("S7").Value = 0.6573033707865168
$ws.Range("T7").Value = 0.6593406593406594
$ws.Range("U7").Value = 0.5508474576271186
$ws.Range("V7").Value = 0.3297101449275363
$ws.Range("W7").Value = 0.4261603375527426
$ws.Range("X7").Value = 0.46
$ws.Range("Y7").Value = 0.68
$ws.Range("Z7").Value = 0.6633663366336634
$ws.Range("AA7").Value = 0.1894736842105263
$ws.Range("AB7").Value = 0.1875
$ws.Range("AD7").Value = 0.703862660944206
$ws.Range("AF7").Value = 0.7951807228915663
$ws.Range("AG7").Value = 0.7840909090909091
$ws.Range("AH7").Value = 0.7727272727272727
$ws.Range("AJ7").Value = 0.8

# Row 8
$ws.Range("B8").Value = 0.8363636363636363
$ws.Range("G8").Value = 0.8906752411575563
$ws.Range("H8").Value = 0.9361702127659575
$ws.Range("I8").Value = 0.9315789473684211
$ws.Range("J8").Value = 0.5312109090909092
$ws.Range("K8").Value = 0.5525508425193906
$ws.Range("L8").Value = 0.8311638711958134
$ws.Range("M8").Value = 0.7997265892002733
$ws.Range("N8").Value = 0.2576923076923077
$ws.Range("O8").Value = 0.006679764243614931
$ws.Range("P8").Value = 0.1642512077294686
$ws.Range("Q8").Value = 0.2315789473684211
$ws.Range("R8").Value = 0.7910447761194029
$ws.Range("S8").Value = 0.7247191011235955
$ws.Range("U8").Value = 0.6440677966101694
$ws.Range("V8").Value = 0.4818840579710146
$ws.Range("W8").Value = 0.5358649789029535
$ws.Range("X8").Value = 0.64
$ws.Range("Y8").Value = 0.8
$ws.Range("Z8").Value = 0.7920792079207921
$ws.Range("AA8").Value = 0.3473684210526315
$ws.Range("AB8").Value = 0.296875
$ws.Range("AD8").Value = 0.7553648068669528
$ws.Range("AE8").Value = 0.8674698795180723
$ws.Range("AF8").Value = 0.855421686746988
$ws.Range("AG8").Value = 0.8522727272727273
$ws.Range("AH8").Value = 0.8295454545454546
$ws.Range("AJ8").Value = 0.8588235294117647

# Row 9
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0.9999999999999999
$ws.Range("AA9").Value = 0.9999999999999999
